$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (want-to-go count) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 769
$wsExhibition.Range("F5").Value = 800
$wsExhibition.Range("F6").Value = 1965
$wsExhibition.Range("F7").Value = 177

# Sheet "全部类型" (sheet4) - update matching rows with the same values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 769
$wsAll.Range("F7").Value = 800
$wsAll.Range("F8").Value = 1965
$wsAll.Range("F10").Value = 177
